# no-op test
$p = $ppt.ActivePresentation
